# New crime data collected — weekly CompStat refresh (022 Pct).
#
# Bumps the report volume/number and the covered week-range in the title
# block, then rolls the Week-to-Date / 28-Day / Year-to-Date / 2-Year crime
# figures (and their dependent % changes) forward for the affected category
# rows. Some cells flip between a numeric count and the sheet's "no data"
# placeholder text ("0" / "***.*" — the same shared text already used
# elsewhere on this sheet, e.g. C14/D14/E14) depending on whether that
# category had any complaints this week. For those, the text is written
# with a leading quote (forces text, not a number) and the number
# format/font is then copied over from a neighboring placeholder cell so it
# keeps matching the rest of the "no activity" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title block: Volume/Number and the covered week date range ----------
# (the "Through" date is edited first so the earlier "Week" date's
# character offset isn't shifted by the length change of 9/8->9/15)
$ws.Range("A8").Characters(21, 2).Text = "38"
$ws.Range("C9").Characters(46, 9).Text = "9/21/2025"
$ws.Range("C9").Characters(27, 8).Text = "9/15/2025"

# --- Row 15 (Rape) --------------------------------------------------------
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("M15").Value2 = -66.666666666666

# --- Row 16 (Robbery) -----------------------------------------------------
# C16/D16/E16 switch from the "no activity" placeholder text back to real
# numbers this week, so restore the usual numeric-column formatting
# (copied from F16/K16, which already carry it) before writing the values.
$ws.Range("F16").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("K16").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C16").Value2 = 1
$ws.Range("D16").Value2 = 1
$ws.Range("E16").Value2 = 0
$ws.Range("F16").Value2 = 2
$ws.Range("H16").Value2 = 100
$ws.Range("I16").Value2 = 7
$ws.Range("J16").Value2 = 35
$ws.Range("K16").Value2 = -80
$ws.Range("L16").Value2 = -58.823529411764
$ws.Range("M16").Value2 = -65
$ws.Range("N16").Value2 = -95.652173913043

# --- Row 17 (Fel. Assault) -------------------------------------------------
$ws.Range("F16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C17").Value2 = 1
$ws.Range("I17").Value2 = 15
$ws.Range("K17").Value2 = 36.363636363636
$ws.Range("L17").Value2 = 114.285714285714
$ws.Range("M17").Value2 = 275
$ws.Range("N17").Value2 = -51.612903225806

# --- Row 18 (Burglary) -----------------------------------------------------
$ws.Range("G18").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("G18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H18").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("H18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("N18").Value2 = -95.652173913043

# --- Row 19 (Gr. Larceny) --------------------------------------------------
$ws.Range("C19").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D19").Value2 = 1
$ws.Range("E19").Value2 = -100
$ws.Range("F19").Value2 = 4
$ws.Range("G19").Value2 = 3
$ws.Range("H19").Value2 = 33.333333333333
$ws.Range("J19").Value2 = 37
$ws.Range("K19").Value2 = -2.702702702702
$ws.Range("N19").Value2 = -74.468085106383

# --- Row 21 (TOTAL) --------------------------------------------------------
$ws.Range("F21").Value2 = 8
$ws.Range("G21").Value2 = 4
$ws.Range("H21").Value2 = 100
$ws.Range("I21").Value2 = 61
$ws.Range("J21").Value2 = 86
$ws.Range("K21").Value2 = -29.069767441860
$ws.Range("L21").Value2 = -7.575757575757
$ws.Range("M21").Value2 = -30.681818181818
$ws.Range("N21").Value2 = -83.646112600536

# --- Row 24 (Petit Larceny) -------------------------------------------------
$ws.Range("C24").Value2 = 2

$ws.Range("D24").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E24").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F24").Value2 = 7
$ws.Range("H24").Value2 = 16.666666666666
$ws.Range("I24").Value2 = 27
$ws.Range("K24").Value2 = -10
$ws.Range("L24").Value2 = -12.903225806451
$ws.Range("M24").Value2 = -57.142857142857

# --- Row 26 (Misd. Assault) -------------------------------------------------
$ws.Range("F16").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C26").Value2 = 1
$ws.Range("D26").Value2 = 1
$ws.Range("E26").Value2 = 0
$ws.Range("F26").Value2 = 1
$ws.Range("G26").Value2 = 4
$ws.Range("H26").Value2 = -75
$ws.Range("I26").Value2 = 22
$ws.Range("J26").Value2 = 30
$ws.Range("K26").Value2 = -26.666666666666
$ws.Range("L26").Value2 = -46.341463414634
$ws.Range("M26").Value2 = 29.411764705882

# --- Row 27 (UCR Rape*) ----------------------------------------------------
$ws.Range("C27").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 28 (Other Sex Crimes) ---------------------------------------------
$ws.Range("C28").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E28").Value2 = -100
$ws.Range("F28").Value2 = 6
$ws.Range("H28").Value2 = 200
$ws.Range("J28").Value2 = 10
$ws.Range("K28").Value2 = 140
